$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Tipo" column (D) to make room for "MAE"
$ws.Range("D1").EntireColumn.Insert()

# New header in the inserted column
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "MAE"

# Update existing values (MSE, R2) and set the new MAE value
$ws.Range("B2").Value = 0.1273602494090013
$ws.Range("C2").Value = 0.9905901608037296
$ws.Range("D2").Value = 0.2777058516343318
